# Updated cryptos list — apply Price (D) and Volume(1h) (E) changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.986.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.93%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.276.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.59%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "506.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.71%  "

$ws.Range("E7").Value = "  -0.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.530"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.289.73"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0988"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.17%  "

$ws.Range("E11").Value = "  +1.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.342"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.681.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "55.017.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.04%  "

$ws.Range("E17").Value = "  +1.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.292.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.39%  "

$ws.Range("E19").Value = "  +1.59%  "

$ws.Range("E20").Value = "  +1.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "314.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "59.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.53%  "

$ws.Range("E26").Value = "  +4.79%  "

$ws.Range("E27").Value = "  +4.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "171.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.68%  "

$ws.Range("E29").Value = "  +4.33%  "

$ws.Range("E30").Value = "  +2.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0707"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.99%  "

$ws.Range("E32").Value = "  +7.68%  "

$ws.Range("E34").Value = "  +1.50%  "

$ws.Range("E35").Value = "  -0.54%  "

$ws.Range("E36").Value = "  +4.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.904"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.81%  "

$ws.Range("E38").Value = "  +5.77%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.49%  "

$ws.Range("E40").Value = "  +4.96%  "

$ws.Range("E41").Value = "  +1.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "136.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.80%  "

$ws.Range("E43").Value = "  +4.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "258.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.13%  "

$ws.Range("E46").Value = "  +3.69%  "

$ws.Range("E47").Value = "  +3.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.549"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.57%  "

$ws.Range("E49").Value = "  +4.58%  "

$ws.Range("E50").Value = "  +1.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.85%  "
